$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# --- "target" column (A): insert "word" between "webcookie" and "ws", shifting the tail down ---
$ws.Range("A30").Value2 = "word"
$ws.Range("A31").Value2 = "ws"
$ws.Range("A32").Value2 = "ws.async"
$ws.Range("A33").Value2 = "xml"

# --- "pdf" column (S): insert "ocr(pdf,saveTo)" between "count(...)" and "saveAsPages(...)", shifting the tail down ---
$ws.Range("S13").Value2 = "ocr(pdf,saveTo)"
$ws.Range("S14").Value2 = "saveAsPages(pdf,destination)"
$ws.Range("S15").Value2 = "saveAsPdf(profile,content,file)"
$ws.Range("S16").Value2 = "saveAsText(pdf,destination)"
$ws.Range("S17").Value2 = "saveFormValues(pdf,var,pageAndLineStartEnd,strategy)"
$ws.Range("S18").Value2 = "saveMetadata(pdf,var)"
$ws.Range("S19").Value2 = "savePageCount(pdf,var)"
$ws.Range("S20").Value2 = "saveToVar(pdf,var)"
$ws.Range("S21").Value2 = "split(pdf,saveTo)"

# --- Shift "xml" list AF->AG, "ws.async" list AE->AF, "ws" list AD->AE (to make room for new "word" list at AD) ---
# column AG
$ws.Range("AG1").Value2 = "xml"
$ws.Range("AG2").Value2 = "append(xml,xpath,content,var)"
$ws.Range("AG3").Value2 = "assertCorrectness(xml,schema)"
$ws.Range("AG4").Value2 = "assertElementCount(xml,xpath,count)"
$ws.Range("AG5").Value2 = "assertElementNotPresent(xml,xpath)"
$ws.Range("AG6").Value2 = "assertElementPresent(xml,xpath)"
$ws.Range("AG7").Value2 = "assertSoap(wsdl,xml)"
$ws.Range("AG8").Value2 = "assertSoapFaultCode(expected,xml)"
$ws.Range("AG9").Value2 = "assertSoapFaultString(expected,xml)"
$ws.Range("AG10").Value2 = "assertValue(xml,xpath,expected)"
$ws.Range("AG11").Value2 = "assertValues(xml,xpath,array,exactOrder)"
$ws.Range("AG12").Value2 = "assertWellformed(xml)"
$ws.Range("AG13").Value2 = "beautify(xml,var)"
$ws.Range("AG14").Value2 = "clear(xml,xpath,var)"
$ws.Range("AG15").Value2 = "delete(xml,xpath,var)"
$ws.Range("AG16").Value2 = "insertAfter(xml,xpath,content,var)"
$ws.Range("AG17").Value2 = "insertBefore(xml,xpath,content,var)"
$ws.Range("AG18").Value2 = "minify(xml,var)"
$ws.Range("AG19").Value2 = "prepend(xml,xpath,content,var)"
$ws.Range("AG20").Value2 = "replace(xml,xpath,content,var)"
$ws.Range("AG21").Value2 = "replaceIn(xml,xpath,content,var)"
$ws.Range("AG22").Value2 = "storeCount(xml,xpath,var)"
$ws.Range("AG23").Value2 = "storeSoapFaultCode(var,xml)"
$ws.Range("AG24").Value2 = "storeSoapFaultDetail(var,xml)"
$ws.Range("AG25").Value2 = "storeSoapFaultString(var,xml)"
$ws.Range("AG26").Value2 = "storeValue(xml,xpath,var)"
$ws.Range("AG27").Value2 = "storeValues(xml,xpath,var)"
# column AF
$ws.Range("AF1").Value2 = "ws.async"
$ws.Range("AF2").Value2 = "delete(url,body,output)"
$ws.Range("AF3").Value2 = "download(url,queryString,saveTo)"
$ws.Range("AF4").Value2 = "get(url,queryString,output)"
$ws.Range("AF5").Value2 = "head(url,output)"
$ws.Range("AF6").Value2 = "patch(url,body,output)"
$ws.Range("AF7").Value2 = "post(url,body,output)"
$ws.Range("AF8").Value2 = "put(url,body,output)"
# column AE
$ws.Range("AE1").Value2 = "ws"
$ws.Range("AE2").Value2 = "assertReturnCode(var,returnCode)"
$ws.Range("AE3").Value2 = "delete(url,body,var)"
$ws.Range("AE4").Value2 = "download(url,queryString,saveTo)"
$ws.Range("AE5").Value2 = "get(url,queryString,var)"
$ws.Range("AE6").Value2 = "head(url,var)"
$ws.Range("AE7").Value2 = "header(name,value)"
$ws.Range("AE8").Value2 = "headerByVar(name,var)"
$ws.Range("AE9").Value2 = "jwtParse(var,token,key)"
$ws.Range("AE10").Value2 = "jwtSignHS256(var,payload,key)"
$ws.Range("AE11").Value2 = "oauth(var,url,auth)"
$ws.Range("AE12").Value2 = "patch(url,body,var)"
$ws.Range("AE13").Value2 = "post(url,body,var)"
$ws.Range("AE14").Value2 = "put(url,body,var)"
$ws.Range("AE15").Value2 = "saveResponsePayload(var,file,append)"
$ws.Range("AE16").Value2 = "soap(action,url,payload,var)"
$ws.Range("AE17").Value2 = "upload(url,body,fileParams,var)"

# --- New "word" list at AD (Word-document automation commands) ---
$ws.Range("AD1").Value2 = "word"
$ws.Range("AD2").Value2 = "assertContains(file,text)"
$ws.Range("AD3").Value2 = "assertNotContain(file,text)"
$ws.Range("AD4").Value2 = "assertNotReadOnly(file)"
$ws.Range("AD5").Value2 = "assertPassword(file,password)"
$ws.Range("AD6").Value2 = "assertReadOnly(file)"
$ws.Range("AD7").Value2 = "extractText(var,file)"
$ws.Range("AD8").Value2 = "readOnly(file,password)"
$ws.Range("AD9").Value2 = "removeProtection(file)"

# --- Update defined names to reflect the new/shifted ranges ---
$wb.Names.Item("pdf").RefersTo = "='#system'!`$S`$2:`$S`$21"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$33"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AE`$2:`$AE`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AF`$2:`$AF`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AG`$2:`$AG`$27"
$wb.Names.Add("word", "='#system'!`$AD`$2:`$AD`$9")
